$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2532
$ws.Range("I62").Value = 2143.5715
$ws.Range("J62").Value = 3075.8
$ws.Range("K62").Value = 2143.5715
$ws.Range("L62").Value = 3075.8
$ws.Range("M62").Value = -1519.5715
$ws.Range("N62").Value = -4323.8

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2532
$ws.Range("I65").Value = 2143.5715
$ws.Range("J65").Value = 3075.8
$ws.Range("K65").Value = 10717.8575
$ws.Range("L65").Value = 15379
$ws.Range("M65").Value = -7597.8575
$ws.Range("N65").Value = -21619

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 5394.4614
$ws.Range("I111").Value = 6260.1
$ws.Range("J111").Value = 2509
$ws.Range("K111").Value = 18780.3
$ws.Range("L111").Value = 7527
$ws.Range("M111").Value = -15713.3
$ws.Range("N111").Value = -13661

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 997.5357
$ws.Range("J112").Value = 1036.1923
$ws.Range("L112").Value = 3108.5769
$ws.Range("N112").Value = -5324.5769

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1666.5
$ws.Range("I137").Value = 2041.5834
$ws.Range("J137").Value = 1216.4
$ws.Range("K137").Value = 6124.7502
$ws.Range("L137").Value = 3649.2
$ws.Range("M137").Value = -3574.7502
$ws.Range("N137").Value = -8749.200000000001

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2290
$ws.Range("I63").Value = 1457.5
$ws.Range("K63").Value = 1457.5
$ws.Range("M63").Value = -771.5

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2290
$ws.Range("I66").Value = 1457.5
$ws.Range("K66").Value = 7287.5
$ws.Range("M66").Value = -3855.5

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1495.8
$ws.Range("I88").Value = 1252.5
$ws.Range("J88").Value = 1658
$ws.Range("K88").Value = 1252.5
$ws.Range("L88").Value = 1658
$ws.Range("M88").Value = -846.5
$ws.Range("N88").Value = -2470

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1495.8
$ws.Range("I91").Value = 1252.5
$ws.Range("J91").Value = 1658
$ws.Range("K91").Value = 1252.5
$ws.Range("L91").Value = 1658
$ws.Range("M91").Value = 151.5
$ws.Range("N91").Value = -4466

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# ARM row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 41255
$ws.Range("J121").Value = 41255
$ws.Range("L121").Value = 41255
$ws.Range("N121").Value = -44749

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2501.621
$ws.Range("I132").Value = 2546.327
$ws.Range("J132").Value = 2335.5715
$ws.Range("K132").Value = 7638.981000000001
$ws.Range("L132").Value = 7006.7145
$ws.Range("M132").Value = -5108.981000000001
$ws.Range("N132").Value = -12066.7145

# BSM row 31
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 30026
$ws.Range("J31").Value = 30026
$ws.Range("L31").Value = 30026
$ws.Range("N31").Value = -30530

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 62812.945
$ws.Range("I86").Value = 80164.71000000001
$ws.Range("K86").Value = 80164.71000000001
$ws.Range("M86").Value = -79041.71000000001

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 62812.945
$ws.Range("I89").Value = 80164.71000000001
$ws.Range("K89").Value = 400823.55
$ws.Range("M89").Value = -395207.55

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 40048.58
$ws.Range("I105").Value = 68382.13
$ws.Range("K105").Value = 68382.13
$ws.Range("M105").Value = -66635.13

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1869.5962
$ws.Range("I134").Value = 1529.5625
$ws.Range("J134").Value = 5950
$ws.Range("K134").Value = 4588.6875
$ws.Range("L134").Value = 17850
$ws.Range("M134").Value = -2053.6875
$ws.Range("N134").Value = -22920

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4549.3335
$ws.Range("I31").Value = 2986.6667
$ws.Range("J31").Value = 4861.8667
$ws.Range("K31").Value = 2986.6667
$ws.Range("L31").Value = 4861.8667
$ws.Range("M31").Value = -2691.6667
$ws.Range("N31").Value = -5451.8667

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4549.3335
$ws.Range("I34").Value = 2986.6667
$ws.Range("J34").Value = 4861.8667
$ws.Range("K34").Value = 2986.6667
$ws.Range("L34").Value = 4861.8667
$ws.Range("M34").Value = -2784.6667
$ws.Range("N34").Value = -5265.8667

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1213.3243
$ws.Range("I5").Value = 1262.0834
$ws.Range("K5").Value = 3786.2502
$ws.Range("M5").Value = -3674.2502

# CUL row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 284.07693
$ws.Range("I14").Value = 284.07693
$ws.Range("K14").Value = 852.2307900000001
$ws.Range("M14").Value = -679.2307900000001

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1213.3243
$ws.Range("I135").Value = 1262.0834
$ws.Range("K135").Value = 11358.7506
$ws.Range("M135").Value = -8823.750599999999

# GSM row 22
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# GSM row 42
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 38761
$ws.Range("J42").Value = 38761
$ws.Range("L42").Value = 38761
$ws.Range("N42").Value = -39731

# GSM row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# GSM row 111
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 70000
$ws.Range("J111").Value = 70000
$ws.Range("L111").Value = 70000
$ws.Range("N111").Value = -76134

# GSM row 112
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 38749
$ws.Range("J112").Value = 38749
$ws.Range("L112").Value = 38749
$ws.Range("N112").Value = -40965

# GSM row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 29722
$ws.Range("J114").Value = 29722
$ws.Range("L114").Value = 29722
$ws.Range("N114").Value = -38400

# GSM row 115
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 38761
$ws.Range("J115").Value = 38761
$ws.Range("L115").Value = 38761
$ws.Range("N115").Value = -41111

# GSM row 117
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 20000
$ws.Range("J117").Value = 20000
$ws.Range("L117").Value = 20000
$ws.Range("N117").Value = -26884

# GSM row 118
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 20500
$ws.Range("J118").Value = 20500
$ws.Range("L118").Value = 20500
$ws.Range("N118").Value = -23814

# GSM row 119
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# GSM row 120
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# GSM row 121
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2231.0435
$ws.Range("I126").Value = 2337.923
$ws.Range("J126").Value = 2092.1
$ws.Range("K126").Value = 7013.768999999999
$ws.Range("L126").Value = 6276.299999999999
$ws.Range("M126").Value = -4543.768999999999
$ws.Range("N126").Value = -11216.3

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2107.5454
$ws.Range("J22").Value = 3000.75
$ws.Range("L22").Value = 3000.75
$ws.Range("N22").Value = -3590.75

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2107.5454
$ws.Range("J27").Value = 3000.75
$ws.Range("L27").Value = 3000.75
$ws.Range("N27").Value = -3214.75

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2452.4
$ws.Range("I136").Value = 2210
$ws.Range("J136").Value = 3422
$ws.Range("K136").Value = 6630
$ws.Range("L136").Value = 10266
$ws.Range("M136").Value = -4080
$ws.Range("N136").Value = -15366

# WVR row 56
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 16601
$ws.Range("J56").Value = 22901.5
$ws.Range("L56").Value = 22901.5
$ws.Range("N56").Value = -24329.5
